# Scheduled-runner update: refresh currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) for a batch of leve rows across several
# class sheets, per the latest market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 618.75
$ws.Range("I12").Value = 437.5
$ws.Range("J12").Value = 800
$ws.Range("K12").Value = 437.5
$ws.Range("L12").Value = 800
$ws.Range("M12").Value = -267.5
$ws.Range("N12").Value = -1140

$ws.Range("H51").Value = 7081.75
$ws.Range("I51").Value = 7001
$ws.Range("J51").Value = 7108.6665
$ws.Range("K51").Value = 7001
$ws.Range("L51").Value = 7108.6665
$ws.Range("M51").Value = -6517
$ws.Range("N51").Value = -8076.6665

$ws.Range("H112").Value = 4012.1052
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 4378.2354
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 13134.7062
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -15350.7062

$ws.Range("H132").Value = 6470.8237
$ws.Range("I132").Value = 5020.614
$ws.Range("J132").Value = 15586.429
$ws.Range("K132").Value = 15061.842
$ws.Range("L132").Value = 46759.287
$ws.Range("M132").Value = -12531.842
$ws.Range("N132").Value = -51819.287

$ws.Range("H138").Value = 2396.51
$ws.Range("I138").Value = 1431.5294
$ws.Range("J138").Value = 2594.1567
$ws.Range("K138").Value = 4294.5882
$ws.Range("L138").Value = 7782.4701
$ws.Range("M138").Value = 845.4117999999999
$ws.Range("N138").Value = -18062.4701

$ws.Range("H141").Value = 1570
$ws.Range("I141").Value = 1570
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4710
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 470
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 11631.904
$ws.Range("I74").Value = 1921.7667
$ws.Range("J74").Value = 35907.25
$ws.Range("K74").Value = 1921.7667
$ws.Range("L74").Value = 35907.25
$ws.Range("M74").Value = -1047.7667
$ws.Range("N74").Value = -37655.25

$ws.Range("H77").Value = 11631.904
$ws.Range("I77").Value = 1921.7667
$ws.Range("J77").Value = 35907.25
$ws.Range("K77").Value = 9608.833499999999
$ws.Range("L77").Value = 179536.25
$ws.Range("M77").Value = -5240.833499999999
$ws.Range("N77").Value = -188272.25

$ws.Range("H132").Value = 1522058.4
$ws.Range("I132").Value = 2226.0193
$ws.Range("J132").Value = 7167150
$ws.Range("K132").Value = 6678.0579
$ws.Range("L132").Value = 21501450
$ws.Range("M132").Value = -4148.0579
$ws.Range("N132").Value = -21506510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12129.187
$ws.Range("I31").Value = 6078.5654
$ws.Range("J31").Value = 19087.4
$ws.Range("K31").Value = 6078.5654
$ws.Range("L31").Value = 19087.4
$ws.Range("M31").Value = -5783.5654
$ws.Range("N31").Value = -19677.4

$ws.Range("H34").Value = 12129.187
$ws.Range("I34").Value = 6078.5654
$ws.Range("J34").Value = 19087.4
$ws.Range("K34").Value = 6078.5654
$ws.Range("L34").Value = 19087.4
$ws.Range("M34").Value = -5876.5654
$ws.Range("N34").Value = -19491.4

$ws.Range("H62").Value = 3529.7273
$ws.Range("I62").Value = 2385.4
$ws.Range("J62").Value = 4483.3335
$ws.Range("K62").Value = 2385.4
$ws.Range("L62").Value = 4483.3335
$ws.Range("M62").Value = -1761.4
$ws.Range("N62").Value = -5731.3335

$ws.Range("H65").Value = 3529.7273
$ws.Range("I65").Value = 2385.4
$ws.Range("J65").Value = 4483.3335
$ws.Range("K65").Value = 11927
$ws.Range("L65").Value = 22416.6675
$ws.Range("M65").Value = -8807
$ws.Range("N65").Value = -28656.6675

$ws.Range("H86").Value = 8090.24
$ws.Range("I86").Value = 9972.75
$ws.Range("J86").Value = 6352.5386
$ws.Range("K86").Value = 9972.75
$ws.Range("L86").Value = 6352.5386
$ws.Range("M86").Value = -8849.75
$ws.Range("N86").Value = -8598.5386

$ws.Range("H89").Value = 8090.24
$ws.Range("I89").Value = 9972.75
$ws.Range("J89").Value = 6352.5386
$ws.Range("K89").Value = 49863.75
$ws.Range("L89").Value = 31762.693
$ws.Range("M89").Value = -44247.75
$ws.Range("N89").Value = -42994.693

$ws.Range("H103").Value = 12888
$ws.Range("I103").Value = 12888
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 12888
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -11716

$ws.Range("H105").Value = 20125.75
$ws.Range("I105").Value = 26251.75
$ws.Range("J105").Value = 13999.75
$ws.Range("K105").Value = 26251.75
$ws.Range("L105").Value = 13999.75
$ws.Range("M105").Value = -24504.75
$ws.Range("N105").Value = -17493.75

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H114").Value = 24626.578
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 24626.578
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 24626.578
$ws.Range("N114").Value = -33304.578

$ws.Range("H118").Value = 70000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 70000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 70000
$ws.Range("N118").Value = -73314

$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H134").Value = 16397230
$ws.Range("I134").Value = 1020.875
$ws.Range("J134").Value = 47628104
$ws.Range("K134").Value = 3062.625
$ws.Range("L134").Value = 142884312
$ws.Range("M134").Value = -527.625
$ws.Range("N134").Value = -142889382

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42876990
$ws.Range("I4").Value = 68469810
$ws.Range("J4").Value = 222294.44
$ws.Range("K4").Value = 205409430
$ws.Range("L4").Value = 666883.3200000001
$ws.Range("M4").Value = -205409318
$ws.Range("N4").Value = -667107.3200000001

$ws.Range("H131").Value = 1468.09
$ws.Range("I131").Value = 762.5
$ws.Range("J131").Value = 1497.4896
$ws.Range("K131").Value = 2287.5
$ws.Range("L131").Value = 4492.468800000001
$ws.Range("M131").Value = 2752.5
$ws.Range("N131").Value = -14572.4688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1031.6
$ws.Range("I9").Value = 1825
$ws.Range("J9").Value = 502.66666
$ws.Range("K9").Value = 1825
$ws.Range("L9").Value = 502.66666
$ws.Range("M9").Value = -1655
$ws.Range("N9").Value = -842.66666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 790308.2
$ws.Range("I132").Value = 1487.6923
$ws.Range("J132").Value = 3353974.8
$ws.Range("K132").Value = 4463.0769
$ws.Range("L132").Value = 10061924.4
$ws.Range("M132").Value = -1933.0769
$ws.Range("N132").Value = -10066984.4

$ws.Range("H136").Value = 11348.591
$ws.Range("I136").Value = 8987.571
$ws.Range("J136").Value = 15480.375
$ws.Range("K136").Value = 26962.713
$ws.Range("L136").Value = 46441.125
$ws.Range("M136").Value = -24412.713
$ws.Range("N136").Value = -46826.685
